$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2097167003757363
$ws.Range("D2").Value = 0.1683836165385131
$ws.Range("E2").Value = 0.1572518775516087
$ws.Range("F2").Value = 1.585520969950892
$ws.Range("G2").Value = 0.9885403027151654
$ws.Range("H2").Value = 0.9588152100885736
$ws.Range("I2").Value = 1.107378677071793
$ws.Range("J2").Value = 0.1853140390249806
$ws.Range("L2").Value = 0.2123279231338344
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("O2").Value = 3.943819789077793

$ws.Range("C3").Value = 0.2089922281124075
$ws.Range("D3").Value = 0.1682660490980155
$ws.Range("E3").Value = 0.1569064231003559
$ws.Range("F3").Value = 1.562029610908397
$ws.Range("G3").Value = 0.9632484458719262
$ws.Range("H3").Value = 0.9514223883851685
$ws.Range("I3").Value = 1.090128063910122
$ws.Range("J3").Value = 0.1846786843179515
$ws.Range("L3").Value = 0.2120131430799503
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("O3").Value = 3.874334649243139

$ws.Range("C4").Value = 0.2086295553525943
$ws.Range("D4").Value = 0.1682508060830088
$ws.Range("E4").Value = 0.156754217974111
$ws.Range("F4").Value = 1.548410802478642
$ws.Range("G4").Value = 0.9482747664937534
$ws.Range("H4").Value = 0.9473179711940674
$ws.Range("I4").Value = 1.080104027687717
$ws.Range("J4").Value = 0.1843640825306565
$ws.Range("L4").Value = 0.2118993815216328
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("O4").Value = 3.83374333843679

$ws.Range("C5").Value = 0.2085024639892339
$ws.Range("D5").Value = 0.1682589471704716
$ws.Range("E5").Value = 0.1567072846195892
$ws.Range("F5").Value = 1.543063623931715
$ws.Range("G5").Value = 0.9423126072119317
$ws.Range("H5").Value = 0.9457548628083714
$ws.Range("I5").Value = 1.076162072721203
$ws.Range("J5").Value = 0.1842548947160054
$ws.Range("L5").Value = 0.2118730520794756
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 3.817723608566013

$ws.Range("C6").Value = 0.2084826120749383
$ws.Range("D6").Value = 0.1682611668560412
$ws.Range("E6").Value = 0.1567004036796718
$ws.Range("F6").Value = 1.542187970625051
$ws.Range("G6").Value = 0.9413310356305118
$ws.Range("H6").Value = 0.9455019260078359
$ws.Range("I6").Value = 1.075516149641338
$ws.Range("J6").Value = 0.1842379134401213
$ws.Range("L6").Value = 0.2118698908420669
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 3.815095055259405

$ws.Range("C7").Value = 0.2086277574812456
$ws.Range("D7").Value = 0.168250857717112
$ws.Range("E7").Value = 0.1567535238725988
$ws.Range("F7").Value = 1.548337868022131
$ws.Range("G7").Value = 0.9481937929152764
$ws.Range("H7").Value = 0.9472964471265044
$ws.Range("I7").Value = 1.080050286271579
$ws.Range("J7").Value = 0.1843625329574934
$ws.Range("L7").Value = 0.2118989452889579
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 3.833525178567527

$ws.Range("C8").Value = 0.2094498752796596
$ws.Range("D8").Value = 0.1683312804555186
$ws.Range("E8").Value = 0.1571203477302276
$ws.Range("F8").Value = 1.577254170375852
$ws.Range("G8").Value = 0.9797043527197076
$ws.Range("H8").Value = 0.956176000301781
$ws.Range("I8").Value = 1.101312856878153
$ws.Range("J8").Value = 0.1850793129123574
$ws.Range("L8").Value = 0.2122029039557418
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("O8").Value = 3.919431255942072

$ws.Range("C9").Value = 0.2117123634842812
$ws.Range("D9").Value = 0.1689394810343856
$ws.Range("E9").Value = 0.1583139770019706
$ws.Range("F9").Value = 1.640343245647259
$ws.Range("G9").Value = 1.0459087253912
$ws.Range("H9").Value = 0.977034632017876
$ws.Range("I9").Value = 1.147513289427891
$ws.Range("J9").Value = 0.1870830904522549
$ws.Range("L9").Value = 0.2134285881979281
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 4.104339501691754

$ws.Range("C10").Value = 0.2137693990294878
$ws.Range("D10").Value = 0.169659388917637
$ws.Range("E10").Value = 0.1594789656587814
$ws.Range("F10").Value = 1.690590273277024
$ws.Range("G10").Value = 1.097250556706456
$ws.Range("H10").Value = 0.9944575747461215
$ws.Range("I10").Value = 1.184206116318762
$ws.Range("J10").Value = 0.1889190421312961
$ws.Range("L10").Value = 0.2147115437645724
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 4.250237925732279

$ws.Range("C11").Value = 0.2147905982925948
$ws.Range("D11").Value = 0.170045870394226
$ws.Range("E11").Value = 0.1600712770171491
$ws.Range("F11").Value = 1.714296040369845
$ws.Range("G11").Value = 1.121196801144606
$ws.Range("H11").Value = 1.00283897150311
$ws.Range("I11").Value = 1.201496678672257
$ws.Range("J11").Value = 0.1898330956426264
$ws.Range("L11").Value = 0.2153779694200111
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 4.318797991820588

$ws.Range("C12").Value = 0.2151895514568452
$ws.Range("D12").Value = 0.1702006735885746
$ws.Range("E12").Value = 0.1603045130335374
$ws.Range("F12").Value = 1.723394706965905
$ws.Range("G12").Value = 1.130349641506712
$ws.Range("H12").Value = 1.006078217826229
$ws.Range("I12").Value = 1.208130248121662
$ws.Range("J12").Value = 0.1901905446574759
$ws.Range("L12").Value = 0.2156422058602843
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 4.345074874831539

$ws.Range("C13").Value = 0.215103085686053
$ws.Range("D13").Value = 0.170166958595324
$ws.Range("E13").Value = 0.1602538842571377
$ws.Range("F13").Value = 1.72142973136927
$ws.Range("G13").Value = 1.128374636983466
$ws.Range("H13").Value = 1.005377681741322
$ws.Range("I13").Value = 1.206697766567885
$ws.Range("J13").Value = 0.1901130586892208
$ws.Range("L13").Value = 0.21558477014716
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 4.339401686917142

$ws.Range("C14").Value = 0.2148231751842218
$ws.Range("D14").Value = 0.1700584369565448
$ws.Range("E14").Value = 0.1600902864486393
$ws.Range("F14").Value = 1.715042152835011
$ws.Range("G14").Value = 1.121948109009395
$ws.Range("H14").Value = 1.003104156372274
$ws.Range("I14").Value = 1.202040703248613
$ws.Range("J14").Value = 0.1898622765234137
$ws.Range("L14").Value = 0.2153994704876681
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 4.320953501781844

$ws.Range("C15").Value = 0.2146533156913222
$ws.Range("D15").Value = 0.1699930639514804
$ws.Range("E15").Value = 0.1599912416341276
$ws.Range("F15").Value = 1.7111454348438
$ws.Range("G15").Value = 1.118022733549765
$ws.Range("H15").Value = 1.001720068927682
$ws.Range("I15").Value = 1.199199315425687
$ws.Range("J15").Value = 0.1897101382579436
$ws.Range("L15").Value = 0.2152875147078461
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 4.309694427688157

$ws.Range("C16").Value = 0.213704377217482
$ws.Range("D16").Value = 0.1696353157803117
$ws.Range("E16").Value = 0.15944150921203
$ws.Range("F16").Value = 1.689058099189793
$ws.Range("G16").Value = 1.095697496723432
$ws.Range("H16").Value = 0.9939189882891242
$ws.Range("I16").Value = 1.183088178493151
$ws.Range("J16").Value = 0.1888608915411112
$ws.Range("L16").Value = 0.2146696547928215
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 4.245801412653179

$ws.Range("C17").Value = 0.2131440901764705
$ws.Range("D17").Value = 0.169430935338184
$ws.Range("E17").Value = 0.1591202178067732
$ws.Range("F17").Value = 1.67572535339248
$ws.Range("G17").Value = 1.08215294468414
$ws.Range("H17").Value = 0.9892498833325192
$ws.Range("I17").Value = 1.173357816738289
$ws.Range("J17").Value = 0.1883600890958448
$ws.Range("L17").Value = 0.2143118030945601
$ws.Range("N17").Value = 3.94211849064385
$ws.Range("O17").Value = 4.207165875450755

$ws.Range("C18").Value = 0.2128298725900777
$ws.Range("D18").Value = 0.1693189355982838
$ws.Range("E18").Value = 0.1589412898942975
$ws.Range("F18").Value = 1.668136550989502
$ws.Range("G18").Value = 1.074418066012612
$ws.Range("H18").Value = 0.9866072293752097
$ws.Range("I18").Value = 1.167817538883241
$ws.Range("J18").Value = 0.1880794658424918
$ws.Range("L18").Value = 0.2141137719158905
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 4.185149898488589

$ws.Range("C19").Value = 0.2127248668036543
$ws.Range("D19").Value = 0.1692819693732446
$ws.Range("E19").Value = 0.1588817170074925
$ws.Range("F19").Value = 1.66558083027995
$ws.Range("G19").Value = 1.071808716799694
$ws.Range("H19").Value = 0.9857198426244622
$ws.Range("I19").Value = 1.165951380452924
$ws.Range("J19").Value = 0.1879857276462076
$ws.Range("L19").Value = 0.214048061832564
$ws.Range("N19").Value = 3.828614786363971
$ws.Range("O19").Value = 4.177731083676065

$ws.Range("C20").Value = 0.2132029013794039
$ws.Range("D20").Value = 0.1694521173488468
$ws.Range("E20").Value = 0.1591538124464265
$ws.Range("F20").Value = 1.677136384886637
$ws.Range("G20").Value = 1.083589032134142
$ws.Range("H20").Value = 0.9897424792849563
$ws.Range("I20").Value = 1.174387797112104
$ws.Range("J20").Value = 0.1884126320830077
$ws.Range("L20").Value = 0.2143490904263388
$ws.Range("N20").Value = 3.957806003281064
$ws.Range("O20").Value = 4.211257359249259

$ws.Range("C21").Value = 0.2149050596126187
$ws.Range("D21").Value = 0.1700900832874837
$ws.Range("E21").Value = 0.1601380966499768
$ws.Range("F21").Value = 1.716915034907316
$ws.Range("G21").Value = 1.123833431493097
$ws.Range("H21").Value = 1.003770172107522
$ws.Range("I21").Value = 1.203406262172692
$ws.Range("J21").Value = 0.1899356303853068
$ws.Range("L21").Value = 0.2154535754410745
$ws.Range("N21").Value = 4.391158149571311
$ws.Range("O21").Value = 4.32636364289317

$ws.Range("C22").Value = 0.2160888860249486
$ws.Range("D22").Value = 0.1705562753564323
$ws.Range("E22").Value = 0.160833483903783
$ws.Range("F22").Value = 1.743622595914971
$ws.Range("G22").Value = 1.150630432038895
$ws.Range("H22").Value = 1.013319185784383
$ws.Range("I22").Value = 1.222872785554273
$ws.Range("J22").Value = 0.190996949417837
$ws.Range("L22").Value = 0.2162446253432933
$ws.Range("N22").Value = 4.673791817957863
$ws.Range("O22").Value = 4.40342629206998

$ws.Range("C23").Value = 0.2154505384675645
$ws.Range("D23").Value = 0.1703029645130556
$ws.Range("E23").Value = 0.1604575838574647
$ws.Range("F23").Value = 1.729303361729549
$ws.Range("G23").Value = 1.136283082525523
$ws.Range("H23").Value = 1.008187866785562
$ws.Range("I23").Value = 1.212437301369036
$ws.Range("J23").Value = 0.1904244769856049
$ws.Range("L23").Value = 0.2158161051854037
$ws.Range("N23").Value = 4.523002190005457
$ws.Range("O23").Value = 4.362128778469753

$ws.Range("C24").Value = 0.2131762882285386
$ws.Range("D24").Value = 0.1694425238239532
$ws.Range("E24").Value = 0.1591386062829621
$ws.Range("F24").Value = 1.676498219639768
$ws.Range("G24").Value = 1.082939614816098
$ws.Range("H24").Value = 0.9895196468177119
$ws.Range("I24").Value = 1.173921975274922
$ws.Range("J24").Value = 0.188388854674578
$ws.Range("L24").Value = 0.2143322088292194
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 4.209406988518879

$ws.Range("C25").Value = 0.2110308521808122
$ws.Range("D25").Value = 0.1687268810625895
$ws.Range("E25").Value = 0.15794040798162
$ws.Range("F25").Value = 1.622592143835831
$ws.Range("G25").Value = 1.027525122433985
$ws.Range("H25").Value = 0.9710231873908697
$ws.Range("I25").Value = 1.134532267181612
$ws.Range("J25").Value = 0.1864770691822457
$ws.Range("L25").Value = 0.2130297470713742
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 4.052554133090439
